# Relationenmodell.docx: "Ergaenze Attribute 'passwort', 'EMail'"
#
# Net changes vs. before.docx:
#   1. The empty paragraph right after the title loses its stray
#      "_GoBack" bookmark (bookmarkStart/bookmarkEnd id=0).
#   2. The "User(...)" paragraph gains ", passwort, EMail" right before
#      the closing ")", and the "_GoBack" bookmark re-appears there
#      (moved, not duplicated).
#   3. The "Device(...)" paragraph loses its stray fr-FR run/paragraph
#      language formatting (picked up by a spell-check pass) and a few
#      of its runs get coalesced.
#   4. The empty paragraph after "Device(...)" loses its fr-FR paragraph
#      mark formatting too.
#
# Each affected paragraph is replaced wholesale (Range.InsertXML) with
# the exact WordprocessingML we want, which keeps this robust to
# Word's run-splitting quirks while still landing on the exact target
# markup.

$d = $word.ActiveDocument

# --- 1. paragraph after the title: drop the "_GoBack" bookmark ---------
$pGoBack = $d.Paragraphs(2)
$pGoBack.Range.InsertXML('<w:p w14:paraId="17C23DCB" w14:textId="11C95BB5" w:rsidR="00A9084D" w:rsidRDefault="00A9084D"><w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p>') | Out-Null

# --- 2. "User(...)" paragraph: add ", passwort, EMail" + bookmark ------
$pUser = $d.Paragraphs(3)
$pUser.Range.InsertXML('<w:p w14:paraId="74CFE2C0" w14:textId="40ACA841" w:rsidR="00A9084D" w:rsidRDefault="00A9084D"><w:r><w:t>User(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00A9084D"><w:rPr><w:u w:val="single"/></w:rPr><w:t>Employee_No</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>First_Name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Last_Name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>passwort</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, EMail</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>)</w:t></w:r></w:p>') | Out-Null

# --- 3. "Device(...)" paragraph: strip fr-FR language formatting -------
$pDevice = $d.Paragraphs(5)
$pDevice.Range.InsertXML('<w:p w14:paraId="65CC28CA" w14:textId="7767970C" w:rsidR="00A9084D" w:rsidRDefault="00A9084D"><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00A9084D"><w:t>Device(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00A9084D"><w:rPr><w:u w:val="single"/></w:rPr><w:t>Id_Device</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00A9084D"><w:t xml:space="preserve">, Type, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00A9084D"><w:t>Last_Maintenance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Next_Maintenance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>&#8593;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Id_Location</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>&#8593;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Employee_No</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>') | Out-Null

# --- 4. empty paragraph after "Device(...)": drop fr-FR paragraph mark -
$pAfterDevice = $d.Paragraphs(6)
$pAfterDevice.Range.InsertXML('<w:p w14:paraId="5C1AA730" w14:textId="058CB1C9" w:rsidR="00C74FAD" w:rsidRDefault="00C74FAD"/>') | Out-Null
